# Update weekly Fruit/Vegetable price data (Perejil, Femacal de La Calera).
# The data rows (2-26) keep their Mercado/Región/Categoría/etc. fixed; only
# the Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M) and Precio $/Kg (P) columns are refreshed
# with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2; D=44165; J=68; K=3000; L=3000; M=3000; P=1000},
    @{Row=3; D=44243; J=45; K=3000; L=3000; M=3000; P=1000},
    @{Row=4; D=44574; J=50; K=3000; L=3000; M=3000; P=1000},
    @{Row=5; D=44193; J=70; K=3000; L=3000; M=3000; P=1000},
    @{Row=6; D=44536; J=125; K=2200; L=2200; M=2200; P=733},
    @{Row=7; D=44559; J=68; K=2000; L=2000; M=2000; P=667},
    @{Row=8; D=44179; J=78; K=3000; L=3000; M=3000; P=1000},
    @{Row=9; D=44224; J=67; K=3000; L=3000; M=3000; P=1000},
    @{Row=10; D=44537; J=88; K=2000; L=2200; M=2091; P=697},
    @{Row=11; D=44225; J=56; K=3000; L=3000; M=3000; P=1000},
    @{Row=12; D=44340; J=54; K=3000; L=3000; M=3000; P=1000},
    @{Row=13; D=44292; J=40; K=3000; L=3000; M=3000; P=1000},
    @{Row=14; D=44242; J=95; K=2500; L=3000; M=2737; P=912},
    @{Row=15; D=44222; J=45; K=3000; L=3000; M=3000; P=1000},
    @{Row=16; D=44627; J=78; K=3500; L=3500; M=3500; P=1167},
    @{Row=17; D=44390; J=50; K=3000; L=3000; M=3000; P=1000},
    @{Row=18; D=44223; J=80; K=2500; L=3000; M=2781; P=927},
    @{Row=19; D=44291; J=45; K=3000; L=3000; M=3000; P=1000},
    @{Row=20; D=44166; J=45; K=2500; L=2500; M=2500; P=833},
    @{Row=21; D=44221; J=50; K=2500; L=2500; M=2500; P=833},
    @{Row=22; D=44389; J=81; K=2800; L=3000; M=2889; P=963},
    @{Row=23; D=44260; J=60; K=3500; L=3500; M=3500; P=1167},
    @{Row=24; D=44669; J=92; K=2500; L=3000; M=2755; P=918},
    @{Row=25; D=44557; J=104; K=2000; L=2500; M=2260; P=753},
    @{Row=26; D=44187; J=65; K=3000; L=3000; M=3000; P=1000}
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value = $r.D   # D: Fecha
    $ws.Cells.Item($row, 10).Value = $r.J  # J: Volumen
    $ws.Cells.Item($row, 11).Value = $r.K  # K: Precio minimo
    $ws.Cells.Item($row, 12).Value = $r.L  # L: Precio maximo
    $ws.Cells.Item($row, 13).Value = $r.M  # M: Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $r.P  # P: Precio $/Kg
}
